# Applies the draft-gandhi-mpls-ioam-sr-06.pptx revision:
#  - Slide 17 ("HbH IOAM Indicator Label - Comparisons"): move the
#    comparison table right and reword the last column header of row 1.
#  - Slide 21 ("Generic PW Control Word ... with IOAM Data Fields"):
#    rename "IOAM Indicator Label" -> "E2E IOAM Indicator Label" in the
#    ASCII packet diagram.
#  - Slide 22 (HbH variant of the same diagram): same rename, plus a
#    narrower diagram box.
#  - Slide 8 ("IOAM Indicator Labels"): bump all body text from 14pt to
#    16pt, which grows the placeholder box to fit.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 17: comparison table
# ---------------------------------------------------------------------
$s17 = $p.Slides.Item(17)
$tbl17Shape = $s17.Shapes.Item(2)

# Nudge the table 0.5" to the right (381000 -> 457200 EMU).
$tbl17Shape.Left = 36.00004

# Row 1 / column 5 header text.
$tbl17 = $tbl17Shape.Table
$cell = $tbl17.Cell(1, 5)
$cell.Shape.TextFrame.TextRange.Text = "Different FIB Entry for Local Label  with and without IOAM Enabled"

# ---------------------------------------------------------------------
# Slide 21: Generic PW Control Word diagram (E2E label rename)
# ---------------------------------------------------------------------
$s21 = $p.Slides.Item(21)
$diagram21 = $s21.Shapes.Item(4)
$tr21 = $diagram21.TextFrame.TextRange
$full21 = $tr21.Text
$oldLine = "   | IOAM Indicator Label                  | TC  |1|  TTL          |"
$newLine = "   | E2E IOAM Indicator Label              | TC  |1|  TTL          |"
$idx21 = $full21.IndexOf($oldLine)
$run21 = $tr21.Characters($idx21 + 1, $oldLine.Length)
$run21.Text = $newLine

# ---------------------------------------------------------------------
# Slide 22: HbH variant diagram (same rename + narrower box)
# ---------------------------------------------------------------------
$s22 = $p.Slides.Item(22)
$diagram22 = $s22.Shapes.Item(4)

# Shrink the diagram box width (5724525 -> 5605463 EMU).
$diagram22.Width = 441.3750763700788

$tr22 = $diagram22.TextFrame.TextRange
$full22 = $tr22.Text
$idx22 = $full22.IndexOf($oldLine)
$run22 = $tr22.Characters($idx22 + 1, $oldLine.Length)
$run22.Text = $newLine

# ---------------------------------------------------------------------
# Slide 8: IOAM Indicator Labels body text 14pt -> 16pt
# ---------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$content8 = $s8.Shapes.Item(2)
$tr8 = $content8.TextFrame.TextRange
$tr8.Font.Size = 16

# The placeholder auto-grows to fit the larger text
# (887129/3437222 -> 887128/3665821 EMU).
$content8.Top = 69.85263942519684
$content8.Height = 288.64736983464564
